$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and G receive numeric-looking text values (e.g. "263.47", "8").
# Excel auto-detects numeric-looking input and stores it as a Number, but the
# source workbook keeps these as Text. Apply a temporary "@" (text) number
# format before writing so the values are stored as text, then clear the
# number format again afterwards so no stray style is left on the cells
# (matches the original un-styled cells).
$numericTextRangeD = $ws.Range("D2:D51")
$numericTextRangeG = $ws.Range("G2:G51")
$numericTextRangeD.NumberFormat = "@"
$numericTextRangeG.NumberFormat = "@"

$ws.Range("D2").Value = '263.47'
$ws.Range("G2").Value = '8'
$ws.Range("D3").Value = '21.18'
$ws.Range("G3").Value = '8'
$ws.Range("D4").Value = '6.115'
$ws.Range("G4").Value = '8'
$ws.Range("D5").Value = '0.06106'
$ws.Range("G5").Value = '8'
$ws.Range("D6").Value = '3.559'
$ws.Range("G6").Value = '8'
$ws.Range("D7").Value = '6.515'
$ws.Range("G7").Value = '8'
$ws.Range("D8").Value = '1.352'
$ws.Range("G8").Value = '8'
$ws.Range("D9").Value = '0.8224'
$ws.Range("G9").Value = '8'
$ws.Range("D10").Value = '0.01330'
$ws.Range("G10").Value = '8'
$ws.Range("D11").Value = '0.1598'
$ws.Range("G11").Value = '8'
$ws.Range("D12").Value = '0.08024'
$ws.Range("G12").Value = '8'
$ws.Range("D13").Value = '0.03426'
$ws.Range("G13").Value = '8'
$ws.Range("D14").Value = '0.03179'
$ws.Range("G14").Value = '8'
$ws.Range("D15").Value = '0.09215'
$ws.Range("G15").Value = '8'
$ws.Range("D16").Value = '3.730'
$ws.Range("G16").Value = '8'
$ws.Range("D17").Value = '0.001635'
$ws.Range("G17").Value = '8'
$ws.Range("D18").Value = '0.04614'
$ws.Range("G18").Value = '8'
$ws.Range("D19").Value = '0.006434'
$ws.Range("G19").Value = '8'
$ws.Range("D20").Value = '0.006141'
$ws.Range("G20").Value = '8'
$ws.Range("G21").Value = '8'
$ws.Range("D22").Value = '0.0001502'
$ws.Range("G22").Value = '8'
$ws.Range("D23").Value = '3.728'
$ws.Range("G23").Value = '8'
$ws.Range("D24").Value = '2.281'
$ws.Range("G24").Value = '8'
$ws.Range("D25").Value = '0.3315'
$ws.Range("G25").Value = '8'
$ws.Range("G26").Value = '8'
$ws.Range("G27").Value = '8'
$ws.Range("G28").Value = '8'
$ws.Range("G29").Value = '8'
$ws.Range("G30").Value = '8'
$ws.Range("G31").Value = '8'
$ws.Range("G32").Value = '8'
$ws.Range("G33").Value = '8'
$ws.Range("G34").Value = '8'
$ws.Range("G35").Value = '8'
$ws.Range("G36").Value = '8'
$ws.Range("G37").Value = '8'
$ws.Range("G38").Value = '8'
$ws.Range("G39").Value = '8'
$ws.Range("D40").Value = '0.04578'
$ws.Range("G40").Value = '8'
$ws.Range("D41").Value = '0.006998'
$ws.Range("G41").Value = '8'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1116'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("G42").Value = '8'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '0.003474'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("G43").Value = '8'
$ws.Range("D44").Value = '0.01053'
$ws.Range("G44").Value = '8'
$ws.Range("D45").Value = '0.00006055'
$ws.Range("G45").Value = '8'
$ws.Range("D46").Value = '0.0009908'
$ws.Range("G46").Value = '8'
$ws.Range("G47").Value = '8'
$ws.Range("D48").Value = '0.8032'
$ws.Range("G48").Value = '8'
$ws.Range("D49").Value = '0.001126'
$ws.Range("G49").Value = '8'
$ws.Range("D50").Value = '0.00001901'
$ws.Range("G50").Value = '8'
$ws.Range("G51").Value = '8'

$numericTextRangeD.ClearFormats()
$numericTextRangeG.ClearFormats()
